# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.898.90'
$ws.Range("E2").Value = '  +1.32%  '

$ws.Range("D3").Value = '1.691.50'
$ws.Range("E3").Value = '  -0.49%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.76%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '315.75'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.39%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.44%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3960'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.00%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3994'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.68%  '

$ws.Range("B9").Value = 'Polygon'
$ws.Range("C9").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.443'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.93%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '52.47'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.28%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.48%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.08728'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.07%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '25.52'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.79%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.401'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.92%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.00001341'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.44%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '7.870'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.08%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.631.63'
$ws.Range("E17").Value = '  -3.74%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '94.90'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.76%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.07208'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.08%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '20.42'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.48%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '7.185'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.64%  '

$ws.Range("E22").Value = '  +0.62%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '14.19'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.09%  '

$ws.Range("D24").Value = '24.915.21'
$ws.Range("E24").Value = '  +1.45%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.395'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.61%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.854'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -5.44%  '

$ws.Range("E27").Value = '  +0.46%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '6.029'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.02%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '162.06'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -3.73%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '148.82'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +3.41%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.062'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.71%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.610'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +20.49%  '

$ws.Range("D33").Value = '1.983.36'
$ws.Range("E33").Value = '  +5.31%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.08514'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.42%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.03111'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.43%  '

$ws.Range("E36").Value = '  -1.98%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '7.047'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.06%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.2869'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.14%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.09682'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +5.25%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '10.79'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.12%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.8097'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -7.58%  '

$ws.Range("E42").Value = '  -2.23%  '

$ws.Range("E43").Value = '  -0.71%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '16.87'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.53%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '2.632'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -1.58%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.7282'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.20%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '4.219'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.95%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.08938'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +8.70%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.389'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.95%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.00%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '139.10'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.27%  '

